$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 125
$ws.Range("H125").Value = 2450
$ws.Range("I125").Value = 3000
$ws.Range("J125").Value = 1900
$ws.Range("K125").Value = 27000
$ws.Range("L125").Value = 17100
$ws.Range("M125").Value = -24540
$ws.Range("N125").Value = -22020
# Row 127
$ws.Range("H127").Value = 1340.6666
$ws.Range("I127").Value = 1209.4
$ws.Range("J127").Value = 1997
$ws.Range("K127").Value = 3628.2
$ws.Range("L127").Value = 5991
$ws.Range("M127").Value = 1331.8
$ws.Range("N127").Value = -15911

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 5000
$ws.Range("I28").Value = 5000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 5000
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -4808
# Row 32
$ws.Range("H32").Value = 1999.5
$ws.Range("I32").Value = 1999.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1999.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1712.5
# Row 45
$ws.Range("H45").Value = 6750
$ws.Range("I45").Value = 2583.3333
$ws.Range("J45").Value = 13000
$ws.Range("K45").Value = 2583.3333
$ws.Range("L45").Value = 13000
$ws.Range("M45").Value = -2206.3333
# Row 50
$ws.Range("H50").Value = 25891.857
$ws.Range("I50").Value = 1624
$ws.Range("J50").Value = 35599
$ws.Range("K50").Value = 1624
$ws.Range("L50").Value = 35599
$ws.Range("M50").Value = -910
$ws.Range("N50").Value = -37027
# Row 99
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 5000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2005
# Row 122
$ws.Range("H122").Value = 4761.5713
$ws.Range("I122").Value = 5055.3335
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 15166.0005
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -12716.0005
# Row 124
$ws.Range("H124").Value = 75799
$ws.Range("I124").Value = 73999
$ws.Range("J124").Value = 76999
$ws.Range("K124").Value = 73999
$ws.Range("L124").Value = 76999
$ws.Range("M124").Value = -69089
$ws.Range("N124").Value = -86819
# Row 125
$ws.Range("H125").Value = 72997
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 72997
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 72997
$ws.Range("M125").Value = -82837

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 978.5714
$ws.Range("I31").Value = 1010.1667
$ws.Range("J31").Value = 789
$ws.Range("K31").Value = 1010.1667
$ws.Range("L31").Value = 789
$ws.Range("M31").Value = -715.1667
# Row 33
$ws.Range("H33").Value = 39162.066
$ws.Range("I33").Value = 9905.166999999999
$ws.Range("J33").Value = 58666.668
$ws.Range("K33").Value = 9905.166999999999
$ws.Range("L33").Value = 58666.668
$ws.Range("M33").Value = -9526.166999999999
# Row 34
$ws.Range("H34").Value = 978.5714
$ws.Range("I34").Value = 1010.1667
$ws.Range("J34").Value = 789
$ws.Range("K34").Value = 1010.1667
$ws.Range("L34").Value = 789
$ws.Range("M34").Value = -808.1667
# Row 58
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
# Row 122
$ws.Range("H122").Value = 719.8
$ws.Range("I122").Value = 719.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2159.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 290.6000000000004
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 2247.5715
$ws.Range("I132").Value = 2247
$ws.Range("J132").Value = 2249.6667
$ws.Range("K132").Value = 6741
$ws.Range("L132").Value = 6749.000100000001
$ws.Range("M132").Value = -4211
$ws.Range("N132").Value = -11809.0001
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 9030
$ws.Range("I3").Value = 9030
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 27090
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -26978
# Row 22
$ws.Range("H22").Value = 3051.7
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3051.7
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 9155.099999999999
$ws.Range("N22").Value = -9493.099999999999
# Row 27
$ws.Range("H27").Value = 3051.7
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3051.7
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 9155.099999999999
$ws.Range("N27").Value = -9359.099999999999
# Row 34
$ws.Range("H34").Value = 4615.143
$ws.Range("I34").Value = 1900
$ws.Range("J34").Value = 5701.2
$ws.Range("K34").Value = 5700
$ws.Range("L34").Value = 17103.6
$ws.Range("M34").Value = -5616
$ws.Range("N34").Value = -17271.6
# Row 39
$ws.Range("H39").Value = 21000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 21000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 63000
$ws.Range("N39").Value = -63588
# Row 55
$ws.Range("H55").Value = 2486.258
$ws.Range("I55").Value = 1330
$ws.Range("J55").Value = 3036.8572
$ws.Range("K55").Value = 3990
$ws.Range("L55").Value = 9110.571599999999
$ws.Range("M55").Value = -3813
$ws.Range("N55").Value = -9464.571599999999
# Row 59
$ws.Range("H59").Value = 1085.7142
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1085.7142
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 3257.1426
$ws.Range("N59").Value = -4337.142599999999
# Row 107
$ws.Range("H107").Value = 542.8
$ws.Range("I107").Value = 79
$ws.Range("J107").Value = 658.75
$ws.Range("K107").Value = 237
$ws.Range("L107").Value = 1976.25
$ws.Range("M107").Value = 1683
$ws.Range("N107").Value = -5816.25
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 10000000
$ws.Range("I10").Value = 10000000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 10000000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -9999831
$ws.Range("N10").ClearContents()
# Row 13
$ws.Range("H13").Value = 124
$ws.Range("I13").Value = 95
$ws.Range("J13").Value = 133.66667
$ws.Range("K13").Value = 95
$ws.Range("L13").Value = 133.66667
$ws.Range("M13").Value = 44
$ws.Range("N13").Value = -411.66667
# Row 43
$ws.Range("H43").Value = 1752
$ws.Range("I43").Value = 1752
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1752
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1601
$ws.Range("N43").ClearContents()
# Row 70
$ws.Range("H70").Value = 125006220
$ws.Range("I70").Value = 8300
$ws.Range("J70").Value = 500000000
$ws.Range("K70").Value = 8300
$ws.Range("L70").Value = 500000000
$ws.Range("M70").Value = -8030
# Row 73
$ws.Range("H73").Value = 125006220
$ws.Range("I73").Value = 8300
$ws.Range("J73").Value = 500000000
$ws.Range("K73").Value = 8300
$ws.Range("L73").Value = 500000000
$ws.Range("M73").Value = -7364
# Row 123
$ws.Range("H123").Value = 63123.25
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 63123.25
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 63123.25
$ws.Range("N123").Value = -68023.25
# Row 134
$ws.Range("H134").Value = 129999
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 129999
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 389997
$ws.Range("N134").Value = -395067

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 470.5
$ws.Range("I55").Value = 460.33334
$ws.Range("J55").Value = 501
$ws.Range("K55").Value = 460.33334
$ws.Range("L55").Value = 501
$ws.Range("M55").Value = -287.33334
$ws.Range("N55").Value = -847
# Row 122
$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5050
$ws.Range("N122").ClearContents()
# Row 127
$ws.Range("H127").Value = 79999.5
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 79999.5
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 79999.5
$ws.Range("N127").Value = -89919.5
# Row 135
$ws.Range("H135").Value = 31000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 31000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 31000
$ws.Range("N135").Value = -41140
# Row 136
$ws.Range("H136").Value = 638790.3
$ws.Range("I136").Value = 456118.53
$ws.Range("J136").Value = 1141137.8
$ws.Range("K136").Value = 1368355.59
$ws.Range("L136").Value = 3423413.4
$ws.Range("M136").Value = -1365805.59
$ws.Range("N136").Value = -3428513.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 6000000
$ws.Range("I3").Value = 6000000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6000000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5999886
# Row 11
$ws.Range("H11").Value = 10069
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 10069
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 10069
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -10353
# Row 131
$ws.Range("H131").Value = 68994
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 68994
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 68994
$ws.Range("N131").Value = -79074
